$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Croatia HNL bases update (11-03-2024 22:32): results/odds recalculated
# for matches 125-128, and two new upcoming fixtures (129, 130) appended.
# Final sheet rows 127..132 fully replace the previous rows 127..129.
# -----------------------------------------------------------------------

# Ensure column A (id) and column E (Date) keep their existing formatting
# (style index 1 = bold/bordered header-like cell, style index 2 = date
# number format) by copying format from an already-styled row instead of
# setting NumberFormat/Font directly (which would mint brand-new style
# entries not present in the target workbook).
foreach ($r in 127..132) {
    $ws.Range("A126").Copy($ws.Range("A$r")) | Out-Null
    $ws.Range("E126").Copy($ws.Range("E$r")) | Out-Null
}

function Set-Row {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

Set-Row 127 @{
    A = 125; B = 6769304; C = "Croatia HNL"; D = "Croatia HNL"; E = 45360.45833333334
    F = "NK Lokomotiva Zagreb"; G = "NK Rudes"
    H = 3; I = 0; J = "H"
    K = 1.333; L = 4.75; M = 7.5
    N = 1.4; O = 4; P = 7.5
    Q = -1.25
    R = 1.975; S = 1.875; T = 2.5; U = 1.95; V = 1.9
    W = 0.3999999999999999; X = -1; Y = -1; Z = 0.9750000000000001; AA = -1
    AB = 0.95; AC = -1
}

Set-Row 128 @{
    A = 126; B = 6788930; C = "Croatia HNL"; D = "Croatia HNL"; E = 45360.54861111111
    F = "HNK Gorica"; G = "Hajduk Split"
    H = 0; I = 3; J = "A"
    K = 6; L = 3.75; M = 1.533
    N = 6.5; O = 3.5; P = 1.55
    Q = 1
    R = 1.825; S = 2.025; T = 2.25; U = 2.025; V = 1.825
    W = -1; X = -1; Y = 0.55; Z = -1; AA = 1.025
    AB = 1.025; AC = -1
}

Set-Row 129 @{
    A = 127; B = 6788932; C = "Croatia HNL"; D = "Croatia HNL"; E = 45361.45833333334
    F = "HNK Rijeka"; G = "NK Osijek"
    H = 3; I = 0; J = "H"
    K = 1.533; L = 3.75; M = 6
    N = 1.571; O = 3.75; P = 5.5
    Q = -0.75
    R = 1.75; S = 2.05; T = 2.5; U = 2; V = 1.85
    W = 0.571; X = -1; Y = -1; Z = 0.75; AA = -1
    AB = 1; AC = -1
}

Set-Row 130 @{
    A = 128; B = 6788931; C = "Croatia HNL"; D = "Croatia HNL"; E = 45361.5625
    F = "Dinamo Zagreb"; G = "Slaven Belupo"
    H = 5; I = 2; J = "H"
    K = 1.25; L = 6; M = 8
    N = 1.166; O = 7; P = 11
    Q = -2
    R = 2; S = 1.85; T = 3; U = 1.925; V = 1.925
    W = 0.1659999999999999; X = -1; Y = -1; Z = 1; AA = -1
    AB = 0.925; AC = -1
}

Set-Row 131 @{
    A = 129; B = 6788933; C = "Croatia HNL"; D = "Croatia HNL"; E = 45366.54166666666
    F = "NK Osijek"; G = "NK Varazdin"
    K = 1.55; L = 3.8; M = 5.5
    N = 1.571; O = 3.8; P = 5.25
    Q = -1
    R = 2.05; S = 1.8; T = 2.5; U = 1.85; V = 2
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}

Set-Row 132 @{
    A = 130; B = 6788935; C = "Croatia HNL"; D = "Croatia HNL"; E = 45368.45833333334
    F = "Istra 1961"; G = "HNK Gorica"
    K = 2; L = 3.2; M = 3.8
    N = 2; O = 3.2; P = 3.8
    Q = -0.5
    R = 2.025; S = 1.825; T = 2.25; U = 1.975; V = 1.875
    W = 0; X = 0; Y = 0; Z = 0; AA = 0
}
